$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Update the email address displayed in column A (rows 2-4) while keeping
# the existing hyperlink targets intact. Only the visible text changes from
# likithal39@gmail.com to likitha.lokesh@slalom.com
$ws.Range("A2").Value = "likitha.lokesh@slalom.com"
$ws.Range("A3").Value = "likitha.lokesh@slalom.com"
$ws.Range("A4").Value = "likitha.lokesh@slalom.com"

# Update the active selection to B5 (was B4)
$ws.Range("B5").Select()

$wb.Save()
